# Pushing the final copy of refactored code
# Adds a new "SignInPage" worksheet (sign-in negative-test data) as the
# last sheet in the workbook, and makes it the active tab.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet after the current last sheet --------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "SignInPage"

# --- Fill in the header row + two data rows --------------------------------
# (write order chosen so the shared-strings table comes out in the same
#  sequence as the authored workbook)
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "expectedMessage"
$ws.Range("D2").Value = "There were errors in your submission`nYour username is a required field`nYour account password is a required field"
$ws.Range("A2").Value = "emptyfields"
$ws.Range("A1").Value = "scenario"
$ws.Range("A3").Value = "invalid"
$ws.Range("B3").Value = "test"
$ws.Range("D3").Value = "There were errors in your submission`nYour username should be a valid email address."
$ws.Range("C3").Value = "test123"

# --- Header row styling: yellow fill, no border -----------------------------
$headerRange = $ws.Range("A1:D1")
$headerRange.Interior.Color = 65535

# --- Body styling: reuse the workbook's existing thin-border cell style ----
# (copy format from a cell that already carries the plain thin border so we
#  don't fork a brand-new, duplicate border definition)
$borderSource = $wb.Worksheets.Item("Tests").Range("A2")
$borderSource.Copy()
$ws.Range("A2:D3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Message columns: wrap text so the multi-line messages are readable ----
$ws.Range("D2").WrapText = $true
$ws.Range("D3").WrapText = $true

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.83
$ws.Columns.Item(2).ColumnWidth = 9
$ws.Columns.Item(3).ColumnWidth = 8.67
$ws.Columns.Item(4).ColumnWidth = 39.67

# --- Row heights for the wrapped message rows -------------------------------
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 45

# --- Leave the selection below the data, matching the authored file --------
$ws.Range("A4").Select()

Write-Output "SignInPage sheet added"
